# Add files via upload
# Fill in the "Ende" (F) column with "Nein" for rows 3 through 26 on Tabelle1,
# and update the active cell selection to F28 (matching the author's last
# click position before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Tabelle1")

for ($r = 3; $r -le 26; $r++) {
    $ws.Cells.Item($r, 6).Value = "Nein"
}

$ws.Activate()
$ws.Range("F28").Select()
